$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows to append (rows 7-9)
$data = @(
    @(9994.06, 9948.2999999999993, 283.47000000000003, 284.77, $false, 0.46, 42613.767106481479, $true),
    @(9996.06, 9994.06, 282.39, 282.45999999999998, $false, 0.02, 42614.674178240741, $true),
    @(10030.049999999999, 9996.06, 280.62, 281.57, $false, 0.34, 42615.752754629626, $true)
)

$r = 7
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item(3, 7).Copy()
    $ws.Cells.Item($r, 7).PasteSpecial(-4122)
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}

$excel.CutCopyMode = $false

